$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.348.87'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.103.33'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '523.82'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.39'
$ws.Range('E6').Value = '  -3.00%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.102.65'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.445'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.27'
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  +2.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.644.16'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('E15').Value = '  -2.38%  '
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.470.98'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.102.21'
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('E19').Value = '  -2.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.38'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.84'
$ws.Range('E21').Value = '  -2.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '346.50'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.56'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0892'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('E31').Value = '  +3.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.86'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.01'
$ws.Range('E33').Value = '  -6.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.84'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.91'
$ws.Range('E35').Value = '  +6.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('E36').Value = '  -3.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.04'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.81'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.16'
$ws.Range('E41').Value = '  +6.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0660'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.60'
$ws.Range('E43').Value = '  +5.26%  '
$ws.Range('E44').Value = '  +2.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.145.29'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '36.54'
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.354.19'
$ws.Range('E47').Value = '  +2.47%  '
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  +3.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.956'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('E51').Value = '  -0.34%  '
